$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The sheet is a task list (A: Task, B: Area, D: Completion Date format
# placeholders). We need to:
#   * rename the "NOT FINISHED" banner text
#   * add 5 new "not finished" tasks at the top of the list (rows 3-7)
#   * add 2 more new tasks interleaved further down
#   * drop the old stray "s" task row
#   * move the selection / window layout
# ---------------------------------------------------------------------

# Step 0: rename the "NOT FINISHED" banner while it is still at its
# original row (row 12); it will be carried down automatically once the
# rows above it are shifted.
$ws.Range("A12").Value = "NOT FINISHED - MORE SOON"

# Step 1: insert 5 new rows at the top of the task list (before the old
# row 3) to hold the new unfinished-work items. Insert() shifts the
# existing formatted rows down and carries row-level formatting with it.
$ws.Rows("3:7").Insert()
$ws.Range("C3:D7").Clear()

$ws.Range("A3").Value = "Fix: Warehouse ramps too steep"
$ws.Range("A4").Value = 'Fix " velocity increasing but not speed" (prediction miss bug when hitting wall at specific angle sometimes)'
$ws.Range("A6").Value = "Complete Release Generation Tool"
$ws.Range("B6").Value = "Engineering"
$ws.Range("A7").Value = "Add kill feed"
$ws.Range("B3").Value = "Bugfix"
$ws.Range("B4").Value = "Bugfix"
$ws.Range("A5").Value = "Add team door brush entity"

$ws.Range("B5").Value = "Feature"
$ws.Range("B7").Value = "Feature"

# Step 2: insert a row for "Finish z_warehouse easter egg" right after
# the z_waves_port row, then insert a row for "Allow people to see what
# team a player is" right before the z_waves_port row.
$ws.Rows("14:14").Insert()
$ws.Range("C14:D14").Clear()

$ws.Rows("13:13").Insert()
$ws.Range("C13:D13").Clear()

$ws.Range("A15").Value = "Finish z_warehouse easter egg"
$ws.Range("A13").Value = "Allow people to see what team a player is"

$ws.Range("B13").Value = "Feature"
$ws.Range("B15").Value = "Content"

# The z_waves_port row (now row 14) no longer carries a completion-date
# placeholder.
$ws.Range("D14").Clear()

# Give the new "Finish z_warehouse easter egg" row the same
# completion-date placeholder formatting (style) as the other date
# placeholder cells, without inventing a new number format.
$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "Properly split out client.h, server.h" row (row 12) gains a
# completion-date placeholder.
$ws.Range("D10").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "z_tdm_warehouse - Second Floor" row (row 8) loses its
# completion-date placeholder.
$ws.Range("D8").Clear()

# Remove the old placeholder sitting just above the "NOT FINISHED" banner
# and collapse the now-empty row so the banner moves up by one row.
$ws.Range("D18").Clear()
$ws.Rows("18:18").Delete()

# Remove the old stray "s" task row entirely.
$ws.Rows("19:19").Delete()

# Drop the left-over completion-date placeholders that used to trail the
# banner; only a single one remains in the final layout.
$ws.Range("D22").Clear()
$ws.Range("D25").Clear()

# ---------------------------------------------------------------------
# View / selection changes
# ---------------------------------------------------------------------
$ws.Range("A13:B15").Select()

$excel.ActiveWindow.WindowState = -4143
$excel.Width = 671
$excel.Height = 637
$excel.Left = 739
$excel.Top = 102
